# ---------------------------------------------------------------------------
# Applies the "final updates to fish survey protocols" edit:
#  1. Body: shorten the "Additional copies..." sentence and split it across
#     five runs (all with identical Arial/22/22 formatting) at the exact
#     text boundaries produced by the authors' edit.
#  2. Body: split "1 person x " into "1 person" + " x " (same text, just
#     two runs instead of one - mirrors a grammar-check style run split).
#  3. Header: rename "Diver Visual Survey" -> "Beach Seines Protocol".
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# Helper-esque pattern: forcing a genuine (but momentary) direct-formatting
# change on a sub-range is the only reliable way to get this COM surface to
# keep two adjacent, identically-formatted runs distinct instead of
# re-coalescing them back into one run.
function Split-Run($startPos, $endPos) {
    $r = $d.Range($startPos, $endPos)
    $r.Bold = 1
    $r.Bold = 0
}

# ---------------------------------------------------------------------------
# 1. "Additional copies of this protocol, field datasheets, data entry
#    templates, instructional videos, literature, and more can be found at:"
#    -> shortened + split into 5 runs.
# ---------------------------------------------------------------------------
$old1 = "Additional copies of this protocol, field datasheets, data entry templates, instructional videos, literature, and more can be found at:"
$seg1a = "Additional copies of this protocol, field datasheets"
$seg1b = " and"
$seg1c = " data entry templates"
$seg1d = " can be found at"
$seg1e = ":"
$new1 = $seg1a + $seg1b + $seg1c + $seg1d + $seg1e

$rng1 = $d.Content
$rng1.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

if ($rng1.Find.Found) {
    $p1Start = $rng1.Start
    $p1End = $rng1.End

    $b1a = $seg1a.Length
    $b1b = $b1a + $seg1b.Length
    $b1c = $b1b + $seg1c.Length
    $b1d = $b1c + $seg1d.Length

    $cut1 = $p1Start + $b1a
    $cut2 = $p1Start + $b1b
    $cut3 = $p1Start + $b1c
    $cut4 = $p1Start + $b1d

    Split-Run $cut1 $cut2
    Split-Run $cut2 $cut3
    Split-Run $cut3 $cut4
    Split-Run $cut4 $p1End
}

# ---------------------------------------------------------------------------
# 2. "1 person x " -> "1 person" + " x " (two runs, same visible text).
# ---------------------------------------------------------------------------
$old2 = "1 person x "
$seg2a = "1 person"

$rng2 = $d.Content
$rng2.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false) | Out-Null

if ($rng2.Find.Found) {
    $p2Start = $rng2.Start
    $p2End = $rng2.End
    $cut5 = $p2Start + $seg2a.Length
    Split-Run $cut5 $p2End
}

# ---------------------------------------------------------------------------
# 3. Header: "Diver Visual Survey" -> "Beach Seines Protocol".
#    Scoped to the primary header range so the similarly-named
#    "Diver Visual Surveys" hyperlink in the document body is untouched.
# ---------------------------------------------------------------------------
$sec = $d.Sections.First
$headers = $sec.Headers
for ($i = 1; $i -le $headers.Count; $i++) {
    $hdr = $headers.Item($i)
    if ($hdr.Exists) {
        $hdr.Range.Find.Execute("Diver Visual Survey", $true, $false, $false, $false, $false, $true, 1, $false, "Beach Seines Protocol", 2) | Out-Null
    }
}
